# Update the "Metadata" worksheet (sheet1)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Version: 5.0.0 -> 6.0.0
$ws.Range("B3").Value = "6.0.0"

# Date: 2021-12-16T17:36:56+00:00 -> 2022-01-21T20:46:54+00:00
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value was empty -> "Alvearie Team"
$ws.Range("B9").Value = "Alvearie Team"

# Row 10 was "Contact" / "No display for ContactDetail" -> "Jurisdiction" / "United States of America"
$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"

# Row 11 was a duplicate "Contact" / "No display for ContactDetail" row - remove it entirely,
# shifting all subsequent rows up by one (21 rows -> 20 rows).
$ws.Rows.Item(11).Delete()

# Update the "Elements" worksheet (sheet2): root Extension row's Short/Definition
# become specific to this profile instead of the generic Extension text.
$ws2 = $wb.Worksheets.Item("Elements")
$ws2.Range("K2").Value = "Procedure Modifier"
$ws2.Range("L2").Value = "Modifier codes used either to supplement information, or to adjust procedure care descriptions"
